{"js": "// Update the worksheet date and all 25 division-problem cells.\n//\n// The table is a 20-row x 5-column grid where only every 4th row (0, 4, 8,\n// 12, 16) actually holds a division problem \u2014 the rows in between are blank\n// spacer rows. We update cells by (row, col) position rather than by\n// matching the old text, because a couple of old values repeat verbatim\n// (e.g. \"54\u00f79=6, 0\" appears twice) but map to different new values\n// depending on which cell they're in.\n\nconst title = { from: \"2024-08-12 Monday\", to: \"2024-08-13 Tuesday\" };\n\n// [tableRowIndex, [newCol0, newCol1, newCol2, newCol3, newCol4]]\nconst rows = [\n  [0, [\"53\u00f75=10, 3\", \"79\u00f75=15, 4\", \"43\u00f76=7, 1\", \"79\u00f75=15, 4\", \"84\u00f78=10, 4\"]],\n  [4, [\"29\u00f77=4, 1\", \"13\u00f74=3, 1\", \"88\u00f77=12, 4\", \"55\u00f74=13, 3\", \"64\u00f77=9, 1\"]],\n  [8, [\"21\u00f75=4, 1\", \"64\u00f73=21, 1\", \"45\u00f79=5, 0\", \"43\u00f79=4, 7\", \"76\u00f76=12, 4\"]],\n  [12, [\"43\u00f79=4, 7\", \"30\u00f76=5, 0\", \"81\u00f79=9, 0\", \"77\u00f76=12, 5\", \"56\u00f77=8, 0\"]],\n  [16, [\"85\u00f73=28, 1\", \"57\u00f74=14, 1\", \"69\u00f77=9, 6\", \"58\u00f74=14, 2\", \"95\u00f74=23, 3\"]],\n];\n\n// Update the title paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nif (titlePara.text === title.from) {\n  titlePara.insertText(title.to, \"Replace\");\n}\n\n// Update the table cells, paragraph by paragraph so run/paragraph\n// formatting (font, size, alignment) is preserved.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (const [rowIndex, newValues] of rows) {\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].insertText(newValues[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all 25 division-problem cells.\n#\n# The table is a 20-row x 5-column grid where only every 4th row (1, 5, 9,\n# 13, 17 in Word's 1-based row numbering) actually holds a division problem\n# -- the rows in between are blank spacer rows. We update cells by (row,\n# col) position rather than by matching the old text, because a couple of\n# old values repeat verbatim (e.g. \"54\u00f79=6, 0\" appears twice) but map to\n# different new values depending on which cell they're in.\n\n$d = $word.ActiveDocument\n\n# Title paragraph (date line).\n$d.Paragraphs.Item(1).Range.Text = \"2024-08-13 Tuesday\"\n\n$table = $d.Tables.Item(1)\n\n# Word-1-based table row => new values for columns 1..5.\n$rowUpdates = @{\n    1  = @(\"53\u00f75=10, 3\", \"79\u00f75=15, 4\", \"43\u00f76=7, 1\", \"79\u00f75=15, 4\", \"84\u00f78=10, 4\")\n    5  = @(\"29\u00f77=4, 1\", \"13\u00f74=3, 1\", \"88\u00f77=12, 4\", \"55\u00f74=13, 3\", \"64\u00f77=9, 1\")\n    9  = @(\"21\u00f75=4, 1\", \"64\u00f73=21, 1\", \"45\u00f79=5, 0\", \"43\u00f79=4, 7\", \"76\u00f76=12, 4\")\n    13 = @(\"43\u00f79=4, 7\", \"30\u00f76=5, 0\", \"81\u00f79=9, 0\", \"77\u00f76=12, 5\", \"56\u00f77=8, 0\")\n    17 = @(\"85\u00f73=28, 1\", \"57\u00f74=14, 1\", \"69\u00f77=9, 6\", \"58\u00f74=14, 2\", \"95\u00f74=23, 3\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $values = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $table.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
